$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row for Rank 8 (META PLATFORMS / FB / 574.41) was accidentally removed.
# It lives at row 9 (row 1 = header, row 2 = rank 1, ... row 9 = rank 8).
$ws.Rows.Item(9).Delete()

# Update the active selection to match the post-edit state.
$ws.Range("H13").Select()
